# "Mejoras de GitHubOrganizationEnquirer y sus tests"
#
# The enquirer test fixture was re-generated by a fresh run of the metrics
# export: the "Metricas guardadas el dia" timestamp in B1, and the matching
# timestamp repeated for every indicator row in column F/G, move forward to
# a newer capture time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("entidadTest")

# Header: "Metricas guardadas el dia" <-> timestamp
$ws.Range("B1").Value = "Fri Mar 08 12:39:48 CET 2024"

# Per-indicator-row timestamp column (F on rows 2-3, G on rows 5-7)
$ws.Range("F2").Value = "Fri Mar 08 11:39:48 CET 2024"
$ws.Range("F3").Value = "Fri Mar 08 11:39:48 CET 2024"
$ws.Range("G5").Value = "Fri Mar 08 11:39:48 CET 2024"
$ws.Range("G6").Value = "Fri Mar 08 11:39:48 CET 2024"
$ws.Range("G7").Value = "Fri Mar 08 11:39:48 CET 2024"
